$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 6527.6665
$ws.Range("J43").Value = 8466.5
$ws.Range("L43").Value = 8466.5
$ws.Range("N43").Value = -8604.5

# Row 70
$ws.Range("H70").Value = 4663.857
$ws.Range("I70").Value = 4384.857
$ws.Range("K70").Value = 13154.571
$ws.Range("M70").Value = -12884.571

# Row 73
$ws.Range("H73").Value = 4663.857
$ws.Range("I73").Value = 4384.857
$ws.Range("K73").Value = 13154.571
$ws.Range("M73").Value = -12218.571

# Row 82
$ws.Range("H82").Value = 46931
$ws.Range("I82").Value = 700
$ws.Range("K82").Value = 2100
$ws.Range("M82").Value = -1694

# Row 85
$ws.Range("H85").Value = 46931
$ws.Range("I85").Value = 700
$ws.Range("K85").Value = 2100
$ws.Range("M85").Value = -696

# Row 98
$ws.Range("H98").Value = 55566164
$ws.Range("I98").Value = 66675396
$ws.Range("J98").Value = 20000
$ws.Range("K98").Value = 66675396
$ws.Range("L98").Value = 20000
$ws.Range("M98").Value = -66673898
$ws.Range("N98").Value = -22996

# Row 122
$ws.Range("H122").Value = 55566164
$ws.Range("I122").Value = 66675396
$ws.Range("J122").Value = 20000
$ws.Range("K122").Value = 200026188
$ws.Range("L122").Value = 60000
$ws.Range("M122").Value = -200023738
$ws.Range("N122").Value = -64900

# Row 125
$ws.Range("H125").Value = 2951.9167
$ws.Range("I125").Value = 2976.5
$ws.Range("K125").Value = 26788.5
$ws.Range("M125").Value = -24328.5

# Row 132
$ws.Range("H132").Value = 4292.8887
$ws.Range("I132").Value = 4082.9285
$ws.Range("J132").Value = 5027.75
$ws.Range("K132").Value = 12248.7855
$ws.Range("L132").Value = 15083.25
$ws.Range("M132").Value = -9718.7855
$ws.Range("N132").Value = -20143.25

# Row 135
$ws.Range("H135").Value = 30416
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 30416
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 273744
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -278814

# Row 138
$ws.Range("H138").Value = 2861.3403
$ws.Range("J138").Value = 3130.4146
$ws.Range("L138").Value = 9391.2438
$ws.Range("N138").Value = -19671.2438

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 13164097
$ws.Range("I32").Value = 13164097
$ws.Range("K32").Value = 13164097
$ws.Range("M32").Value = -13163810

# Row 94
$ws.Range("H94").Value = 46744.5
$ws.Range("J94").Value = 46744.5
$ws.Range("L94").Value = 46744.5
$ws.Range("N94").Value = -48546.5

# Row 132
$ws.Range("H132").Value = 4725.59
$ws.Range("I132").Value = 2870.5667
$ws.Range("J132").Value = 10909
$ws.Range("K132").Value = 8611.7001
$ws.Range("L132").Value = 32727
$ws.Range("M132").Value = -6081.7001
$ws.Range("N132").Value = -37787

$ws = $wb.Worksheets.Item("BSM")
# Row 100
$ws.Range("H100").Value = 14900
$ws.Range("J100").Value = 14900
$ws.Range("L100").Value = 14900
$ws.Range("N100").Value = -17064

# Row 134
$ws.Range("H134").Value = 82408
$ws.Range("I134").Value = 2129
$ws.Range("J134").Value = 350004.66
$ws.Range("K134").Value = 6387
$ws.Range("L134").Value = 1050013.98
$ws.Range("M134").Value = -3852
$ws.Range("N134").Value = -1055083.98

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2537.7942
$ws.Range("I58").Value = 2435.8462
$ws.Range("K58").Value = 2435.8462
$ws.Range("M58").Value = -2232.8462

# Row 105
$ws.Range("H105").Value = 1742.5
$ws.Range("I105").Value = 1406.4286
$ws.Range("K105").Value = 1406.4286
$ws.Range("M105").Value = 340.5714

# Row 132
$ws.Range("H132").Value = 2552.9375
$ws.Range("I132").Value = 2552.9375
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7658.8125
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5128.8125
$ws.Range("N132").ClearContents()

# Row 134
$ws.Range("H134").Value = 459038.88
$ws.Range("I134").Value = 626871.5600000001
$ws.Range("K134").Value = 1880614.68
$ws.Range("M134").Value = -1878079.68

# Row 136
$ws.Range("H136").Value = 2537.7942
$ws.Range("I136").Value = 2435.8462
$ws.Range("K136").Value = 7307.5386
$ws.Range("M136").Value = -4757.5386

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 16126.16
$ws.Range("J2").Value = 28711.357
$ws.Range("L2").Value = 172268.142
$ws.Range("N2").Value = -172494.142

# Row 12
$ws.Range("H12").Value = 634609.25
$ws.Range("J12").Value = 950706.2
$ws.Range("L12").Value = 2852118.6
$ws.Range("N12").Value = -2852464.6

# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

# Row 97
$ws.Range("H97").Value = 1390.8182
$ws.Range("I97").Value = 942
$ws.Range("J97").Value = 1647.2858
$ws.Range("K97").Value = 2826
$ws.Range("L97").Value = 4941.857400000001
$ws.Range("M97").Value = -2330
$ws.Range("N97").Value = -5933.857400000001

# Row 132
$ws.Range("H132").Value = 2259.7144
$ws.Range("J132").Value = 2067
$ws.Range("L132").Value = 18603
$ws.Range("N132").Value = -23663

# Row 137
$ws.Range("H137").Value = 5500.5
$ws.Range("J137").Value = 5222.875
$ws.Range("L137").Value = 15668.625
$ws.Range("N137").Value = -25868.625

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 1528.9584
$ws.Range("J122").Value = 2250
$ws.Range("L122").Value = 6750
$ws.Range("N122").Value = -11650

# Row 132
$ws.Range("H132").Value = 71431930
$ws.Range("I132").Value = 111114776
$ws.Range("J132").Value = 2792.6
$ws.Range("K132").Value = 333344328
$ws.Range("L132").Value = 8377.799999999999
$ws.Range("M132").Value = -333341798
$ws.Range("N132").Value = -13437.8

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 15466500
$ws.Range("I7").Value = 25003872
$ws.Range("J7").Value = 206703.8
$ws.Range("K7").Value = 25003872
$ws.Range("L7").Value = 206703.8
$ws.Range("M7").Value = -25003760
$ws.Range("N7").Value = -206927.8

# Row 22
$ws.Range("H22").Value = 3214.4285
$ws.Range("I22").Value = 3950.4
$ws.Range("K22").Value = 3950.4
$ws.Range("M22").Value = -3655.4

# Row 27
$ws.Range("H27").Value = 3214.4285
$ws.Range("I27").Value = 3950.4
$ws.Range("K27").Value = 3950.4
$ws.Range("M27").Value = -3843.4

# Row 46
$ws.Range("H46").Value = 3404.95
$ws.Range("J46").Value = 3224.75
$ws.Range("L46").Value = 3224.75
$ws.Range("N46").Value = -3600.75

# Row 75
$ws.Range("H75").Value = 86500
$ws.Range("I75").Value = 60000
$ws.Range("K75").Value = 60000
$ws.Range("M75").Value = -59064

# Row 78
$ws.Range("H78").Value = 86500
$ws.Range("I78").Value = 60000
$ws.Range("K78").Value = 180000
$ws.Range("M78").Value = -175320

# Row 126
$ws.Range("H126").Value = 15466500
$ws.Range("I126").Value = 25003872
$ws.Range("J126").Value = 206703.8
$ws.Range("K126").Value = 75011616
$ws.Range("L126").Value = 620111.3999999999
$ws.Range("M126").Value = -75009146
$ws.Range("N126").Value = -625051.3999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 38
$ws.Range("H38").Value = 10062
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 10062
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 10062
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -11008

# Row 124
$ws.Range("H124").Value = 82791.8
$ws.Range("J124").Value = 82791.8
$ws.Range("L124").Value = 82791.8
$ws.Range("N124").Value = -92611.8

# Row 126
$ws.Range("H126").Value = 1533.5
$ws.Range("I126").Value = 1533.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4600.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2130.5
$ws.Range("N126").ClearContents()

# Row 132
$ws.Range("H132").Value = 1627.3
$ws.Range("I132").Value = 1697.3334
$ws.Range("J132").Value = 997
$ws.Range("K132").Value = 5092.0002
$ws.Range("L132").Value = 2991
$ws.Range("M132").Value = -2562.0002
$ws.Range("N132").Value = -8051
